$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.954273223876953
$ws.Range("B1").Value = 2.416875600814819
$ws.Range("C1").Value = 3.361053943634033
$ws.Range("D1").Value = 6.471306324005127
$ws.Range("E1").Value = 1.749944925308228
